$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert: remove the most recently added meeting row (row 2, "第596回") and
# shift everything else up, restoring the previous state where "第595回" was
# the newest entry. At that point in time, the minutes ("議事録") for the two
# newest meetings ("第595回" and "第594回", now rows 2 and 3) had not yet been
# published, so column D ("議事録／議事要旨") reverts to "－" for those rows.

$ws.Rows(2).Delete()

$ws.Range("D2").Value = "－"
$ws.Range("D3").Value = "－"
